$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rows 4-6: round the "Ost" (Q) / "Nord" (R) coordinates to whole numbers,
# and remove the "Starttid" (Z) / "Sluttid" (AB) time values entirely.
foreach ($r in 4..6) {
    $ws.Range("Q$r").Value = [Math]::Round($ws.Range("Q$r").Value(), 0)
    $ws.Range("R$r").Value = [Math]::Round($ws.Range("R$r").Value(), 0)
    $ws.Range("Z$r").ClearContents()
    $ws.Range("AB$r").ClearContents()
}
